# Fix the capitalisation of the Abstract's opening words ("STORIES ABOUT" ->
# "Stories about"). Word's hidden "_GoBack" bookmark always tracks the most
# recent edit location, so once we touch the text, the bookmark relocates
# from the empty paragraph above (right after the author table) to sit right
# after the newly-edited words, splitting the run in two.
#
# We reproduce that by first (re)placing the "_GoBack" bookmark at the edit
# point -- which both moves it off its old spot and forces a run split there
# -- and only then rewriting the text of the leading words, so the run
# carrying the remainder of the sentence stays untouched.

$d = $word.ActiveDocument

$findRange = $d.Content
$findRange.Find.Execute("STORIES ABOUT", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "", 0) | Out-Null

$splitPoint = $findRange.End
$bookmarkRange = $d.Range($splitPoint, $splitPoint)
$d.Bookmarks.Add("_GoBack", $bookmarkRange) | Out-Null

$wordsRange = $d.Content
$wordsRange.Find.Execute("STORIES ABOUT", $true, $false, $false, $false, $false, `
                          $true, 1, $false, "", 0) | Out-Null
$wordsRange.Text = "Stories about"
